$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: merge a name that was split across multiple runs back into a single
# run by doing a literal Find/Replace over the full (already-correct) text.
# ---------------------------------------------------------------------------
function Merge-Name($fullName) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($fullName, $true, $false, $false, $false, $false, $true, 1, $false, $fullName, 2) | Out-Null
}

Merge-Name("Kareem Sherif El-Meteny")
Merge-Name("Yassin Amr El-Helly")

# ---------------------------------------------------------------------------
# Table 1 ("Leader"): resize the first three columns, append a fourth
# "Component" column and fill in its values.
# ---------------------------------------------------------------------------
$t1 = $d.Tables.Item(1)

$t1.Columns(1).Width = 3532 / 20
$t1.Columns(2).Width = 1211 / 20
$t1.Columns(3).Width = 2415 / 20

$t1.Columns.Add() | Out-Null
$t1.Columns(4).Width = 2192 / 20

$t1Header = $t1.Cell(1, 4)
$t1Header.Range.Text = "Component"
$t1Header.Range.Font.Bold = $true

$t1.Cell(2, 4).Range.Text = "Student 2"

# ---------------------------------------------------------------------------
# Table 2 ("Members"): same treatment, five rows of data.
# ---------------------------------------------------------------------------
$t2 = $d.Tables.Item(2)

$t2.Columns(1).Width = 3580 / 20
$t2.Columns(2).Width = 1204 / 20
$t2.Columns(3).Width = 2397 / 20

$t2.Columns.Add() | Out-Null
$t2.Columns(4).Width = 2169 / 20

$t2Header = $t2.Cell(1, 4)
$t2Header.Range.Text = "Component"
$t2Header.Range.Font.Bold = $true

$t2.Cell(2, 4).Range.Text = "Student 1"
$t2.Cell(3, 4).Range.Text = "Admin 1"
$t2.Cell(4, 4).Range.Text = "Admin 2"
$t2.Cell(5, 4).Range.Text = "Advisor"

Write-Output "done"
